# Apply "lag price and region fixed effect" restructuring to Sheet1.
#
# The sheet is rebuilt from scratch: a new row 1 adds three merged,
# centered block headers ("R squared" / "Predicting insecure" /
# "Overall accuracy ") above the existing logFCS/HDDS/rCSI column labels,
# which are now repeated across all three blocks (cols C:E, G:I, K:M).
# Two new "year.ols" variant rows are inserted (lag price, GIEWS price x2)
# plus a new "year.ols+quarter*region FE" row, and the trailing random.*
# rows / year.LASSO row are pushed far down the sheet (rows 25, 32-35),
# matching the final row layout produced by the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - the row/column layout changes too much to edit in place.
$ws.Cells.Clear()

# ---- Row 1: merged 3-block headers ----
$ws.Range("C1").Value = "R squared"
$ws.Range("G1").Value = "Predicting insecure"
$ws.Range("K1").Value = "Overall accuracy "

# Apply the centered alignment to every individual cell first so the merged
# blocks all resolve to the same cell-format record, then merge.
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("L1").HorizontalAlignment = -4108
$ws.Range("M1").HorizontalAlignment = -4108

$ws.Range("C1:E1").Merge()
$ws.Range("G1:I1").Merge()
$ws.Range("K1:M1").Merge()

# ---- Row 2: column labels, repeated across the three blocks ----
$ws.Range("C2").Value = "logFCS"
$ws.Range("D2").Value = "HDDS"
$ws.Range("E2").Value = "rCSI"
$ws.Range("G2").Value = "logFCS"
$ws.Range("H2").Value = "HDDS"
$ws.Range("I2").Value = "rCSI"
$ws.Range("K2").Value = "logFCS"
$ws.Range("L2").Value = "HDDS"
$ws.Range("M2").Value = "rCSI"

# ---- Row 3 ----
$ws.Range("B3").Value = "Old result "
$ws.Range("C3").Value = 0.536
$ws.Range("D3").Value = 0.623
$ws.Range("E3").Value = 0.169

# ---- Row 5 ----
$ws.Range("B5").Value = "year.ols"
$ws.Range("C5").Value = 0.5487049
$ws.Range("D5").Value = 0.5688843
$ws.Range("E5").Value = 0.1299255
$ws.Range("C5:E5").NumberFormat = "0.000"

# ---- Row 6: new - year.ols + lag price (no data yet) ----
$ws.Range("B6").Value = "year.ols + lag price"
$ws.Range("C6:E6").NumberFormat = "0.000"

# ---- Row 7: new - year.ols + GIEWS price (no data yet) ----
$ws.Range("B7").Value = "year.ols + GIEWS price"
$ws.Range("C7:E7").NumberFormat = "0.000"

# ---- Row 8: new - year.ols + GIEWS price (no data yet) ----
$ws.Range("B8").Value = "year.ols + GIEWS price"
$ws.Range("C8:E8").NumberFormat = "0.000"

# ---- Row 9 ----
$ws.Range("B9").Value = "year.ols+ quarterFE"
$ws.Range("C9").Value = 0.559955
$ws.Range("D9").Value = 0.6008912
$ws.Range("E9").Value = 0.1272171
$ws.Range("C9:E9").NumberFormat = "0.000"

# ---- Row 10 ----
$ws.Range("B10").Value = "year.ols+ monthFE "
$ws.Range("C10").Value = 0.56921231
$ws.Range("D10").Value = 0.5951778
$ws.Range("E10").Value = 0.1001644
$ws.Range("C10:E10").NumberFormat = "0.000"

# ---- Row 11: new header-only label, no data cells ----
$ws.Range("B11").Value = "year.ols+quarter*region FE"

# ---- Rows 12-13: styled but empty placeholder rows ----
$ws.Range("C12:E13").NumberFormat = "0.000"

# ---- Row 25 ----
$ws.Range("B25").Value = "year.LASSO"
$ws.Range("C25").Value = 0.56968424
$ws.Range("D25").Value = 0.6002007
$ws.Range("E25").Value = 0.11101
$ws.Range("C25:E25").NumberFormat = "0.000"

# ---- Row 32 ----
$ws.Range("B32").Value = "random.ols"
$ws.Range("C32").Value = 0.4970718
$ws.Range("D32").Value = 0.5477649
$ws.Range("E32").Value = 0.0950986
$ws.Range("C32:E32").NumberFormat = "0.000"

# ---- Row 33 ----
$ws.Range("B33").Value = "random.ols+ quarterFE"
$ws.Range("C33").Value = 0.5181456
$ws.Range("D33").Value = 0.5407152
$ws.Range("E33").Value = 0.1117934
$ws.Range("C33:E33").NumberFormat = "0.000"

# ---- Row 34 ----
$ws.Range("B34").Value = "random.ols+ monthFE+quaterFE"
$ws.Range("C34").Value = 0.526206
$ws.Range("D34").Value = 0.5419876
$ws.Range("E34").Value = 0.1802338
$ws.Range("C34:E34").NumberFormat = "0.000"

# ---- Row 35 ----
$ws.Range("B35").Value = "random.LASSO "
$ws.Range("C35").Value = 0.5078364
$ws.Range("D35").Value = 0.5592047
$ws.Range("E35").Value = 0.1662719
$ws.Range("C35:E35").NumberFormat = "0.000"

# ---- Selection cursor moved as part of the edit ----
$ws.Range("F18").Select()
